# Ecken werden in der oberen Haelfte des Rechtecks detektiert
# Insert a new "Anzahl Ecken" column before the existing "Labels" column (D -> E),
# and populate it with per-row corner counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "Labels" column (D) one to the right, making room for
# the new "Anzahl Ecken" column at D.
$ws.Columns.Item(4).Insert()

# New header cell for the inserted column, matching the style of the other
# header cells (bold / centered / bordered), same as the old D1 style.
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Cells.Item(1, 4).Value = "Anzahl Ecken"

# Per-row corner counts for rows 2..100 (column D), numeric values.
$ankenWerte = @(4,3,2,6,3,4,3,4,5,1,6,1,6,8,3,3,4,3,3,3,4,1,5,2,7,11,3,8,4,4,9,8,7,6,1,3,4,7,2,5,3,1,5,2,4,4,7,1,6,1,1,4,5,3,3,1,2,6,5,4,6,8,5,3,8,4,2,1,1,5,9,2,1,1,2,2,6,2,0,0,1,1,5,0,1,7,2,2,1,6,2,2,2,4,1,1,3,1,4)

for ($i = 0; $i -lt $ankenWerte.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $ankenWerte[$i]
}

$ws.Range("A1").Select()
